# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data table (rows 88-89),
# pushing the existing Femacal de La Calera - Papaya records down by two
# rows (old row 88 -> new row 90, ... old row 125 -> new row 127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 88..125 down to 90..127, leaving two blank rows (88:89) for the
# new records.
$ws.Rows("88:89").Insert()

# New row 88 - Primera
$ws.Cells.Item(88, 1).Value = 3
$ws.Cells.Item(88, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 45215
$ws.Cells.Item(88, 5).Value = 5
$ws.Cells.Item(88, 6).Value = "Fruta"
$ws.Cells.Item(88, 7).Value = 100108
$ws.Cells.Item(88, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(88, 9).Value = 100108004
$ws.Cells.Item(88, 10).Value = "Papaya"
$ws.Cells.Item(88, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(88, 12).Value = "Primera"
$ws.Cells.Item(88, 13).Value = 53
$ws.Cells.Item(88, 14).Value = 16000
$ws.Cells.Item(88, 15).Value = 16000
$ws.Cells.Item(88, 16).Value = 16000
$ws.Cells.Item(88, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(88, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(88, 19).Value = 1600
$ws.Cells.Item(88, 20).Value = 10

# New row 89 - Segunda
$ws.Cells.Item(89, 1).Value = 3
$ws.Cells.Item(89, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(89, 3).Value = "Coquimbo"
$ws.Cells.Item(89, 4).Value = 45215
$ws.Cells.Item(89, 5).Value = 5
$ws.Cells.Item(89, 6).Value = "Fruta"
$ws.Cells.Item(89, 7).Value = 100108
$ws.Cells.Item(89, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(89, 9).Value = 100108004
$ws.Cells.Item(89, 10).Value = "Papaya"
$ws.Cells.Item(89, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(89, 12).Value = "Segunda"
$ws.Cells.Item(89, 13).Value = 45
$ws.Cells.Item(89, 14).Value = 13000
$ws.Cells.Item(89, 15).Value = 13000
$ws.Cells.Item(89, 16).Value = 13000
$ws.Cells.Item(89, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(89, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(89, 19).Value = 1300
$ws.Cells.Item(89, 20).Value = 10
